# Bond dates update: one additional day has elapsed since the figures were
# last computed, so:
#   - "Dni od poprzedniej wypłaty" (days since previous payout, column G)
#     increases by 1 for every bond that has a previous-payout date.
#   - "Dni do następnej wypłaty" (days until next payout, column I)
#     decreases by 1 for every bond that has a next-payout date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $iCell = $ws.Cells.Item($r, 9)   # column I

    $gVal = $gCell.Value2
    $iVal = $iCell.Value2

    if ($gVal -ne $null) {
        $gCell.Value2 = $gVal + 1
    }

    if ($iVal -ne $null) {
        $iCell.Value2 = $iVal - 1
    }
}
